# Horarios actualizados Linea 141 - 822
# Applies the scraped-data refresh (new "Ultima actualizacion" timestamp,
# updated Minutos values for existing rows, and newly-scraped rows appended
# in Hora_Llegada order) across the three schedule sheets.

$wb = $excel.ActiveWorkbook


# ---- Sheet "LP1912" ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 06:21:22"
$ws1.Range("A3").Value = "Total filas: 64"

$sheet1Rows = @(
  @('06:21:22', '06:21', '26_HERNANDEZ', 0, 'LP1912'),
  @('04:48:57', '06:26', '23_HERNANDEZ', 98, 'LP1912'),
  @('05:52:07', '06:27', '23_HERNANDEZ', 35, 'LP1912'),
  @('06:21:22', '06:29', '86_EST CHICA-ESC AGRARIA', 8, 'LP1912'),
  @('06:21:22', '06:29', '23_HERNANDEZ', 8, 'LP1912'),
  @('05:52:07', '06:30', '86_EST CHICA-ESC AGRARIA', 38, 'LP1912'),
  @('06:21:22', '06:31', '16_SANTA ANA', 10, 'LP1912'),
  @('04:48:57', '06:43', '225_C ROCA-H SUR', 115, 'LP1912'),
  @('06:21:22', '06:44', '225_C ROCA-H SUR', 23, 'LP1912'),
  @('06:21:22', '06:46', '215C_EL PATO', 25, 'LP1912'),
  @('05:52:07', '06:47', '215C_EL PATO', 55, 'LP1912'),
  @('06:21:22', '06:59', '14_ABASTO', 38, 'LP1912'),
  @('05:52:07', '07:00', '14_ABASTO', 68, 'LP1912'),
  @('06:21:22', '07:01', '16_SANTA ANA', 40, 'LP1912'),
  @('05:52:07', '07:05', '23_HERNANDEZ', 73, 'LP1912'),
  @('06:21:22', '07:05', '15_ABASTO', 44, 'LP1912'),
  @('06:21:22', '07:07', '225_GOMEZ', 46, 'LP1912'),
  @('06:21:22', '07:11', '215A_EL PATO', 50, 'LP1912'),
  @('05:52:07', '07:12', '215A_EL PATO', 80, 'LP1912'),
  @('06:21:22', '07:15', '11_ETCHEVERRY', 54, 'LP1912'),
  @('05:52:07', '07:16', '11_ETCHEVERRY', 84, 'LP1912'),
  @('06:21:22', '07:21', '26_HERNANDEZ', 60, 'LP1912'),
  @('06:21:22', '07:23', '10_OLMOS', 62, 'LP1912'),
  @('06:21:22', '07:31', '11_ETCHEVERRY', 70, 'LP1912'),
  @('05:52:07', '07:32', '16_SANTA ANA', 100, 'LP1912'),
  @('06:21:22', '07:32', '84_COLONIA URQUIZA-ESC 49', 71, 'LP1912'),
  @('05:52:07', '07:32', '11_ETCHEVERRY', 100, 'LP1912'),
  @('06:21:22', '07:37', '27_EL RETIRO', 76, 'LP1912'),
  @('06:21:22', '07:39', '10_OLMOS', 78, 'LP1912'),
  @('06:21:22', '07:47', '14_ABASTO', 86, 'LP1912'),
  @('05:52:07', '07:48', '14_ABASTO', 116, 'LP1912'),
  @('06:21:22', '07:51', '215D_EL PATO', 90, 'LP1912'),
  @('06:21:22', '08:01', '23_HERNANDEZ', 100, 'LP1912'),
  @('06:21:22', '08:12', '15_ABASTO', 111, 'LP1912')
)

$startRow1 = 36
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $r = $startRow1 + $i
    $row = $sheet1Rows[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}


# ---- Sheet "LP1912-215" ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 06:21:22"
$ws2.Range("A3").Value = "Total filas: 15"

$ws2.Cells.Item(16, 1).Value = '06:21:22'
$ws2.Cells.Item(16, 2).Value = '06:46'
$ws2.Cells.Item(16, 3).Value = '215C_EL PATO'
$ws2.Cells.Item(16, 4).Value = 25
$ws2.Cells.Item(16, 5).Value = 'LP1912'
$ws2.Cells.Item(18, 1).Value = '06:21:22'
$ws2.Cells.Item(18, 2).Value = '07:11'
$ws2.Cells.Item(18, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(18, 4).Value = 50
$ws2.Cells.Item(18, 5).Value = 'LP1912'
$ws2.Cells.Item(20, 1).Value = '06:21:22'
$ws2.Cells.Item(20, 2).Value = '07:51'
$ws2.Cells.Item(20, 3).Value = '215D_EL PATO'
$ws2.Cells.Item(20, 4).Value = 90
$ws2.Cells.Item(20, 5).Value = 'LP1912'

# ---- Sheet "6203-6173" ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:21:22"
$ws3.Range("A3").Value = "Total filas: 11"

$sheet3Rows13to16 = @(
  @('06:21:22', '06:33', '215C_LA PLATA', 12, 'L6203'),
  @('06:21:22', '07:00', '215B_LP-P MOR-1 Y 57', 39, 'L6173'),
  @('06:21:22', '07:35', '215A_LA PLATA', 74, 'L6173'),
  @('06:21:22', '08:07', '215C_LA PLATA', 106, 'L6203')
)

$startRow3 = 13
for ($i = 0; $i -lt $sheet3Rows13to16.Count; $i++) {
    $r = $startRow3 + $i
    $row = $sheet3Rows13to16[$i]
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
}
